$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    if ($r.Text -like "*Lire beaucoup des livres*") {
        $r.Font.Size = 9999999
    }
}
